$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 338.6742503333333
$ws.Range("H2").Value = 1016.022751
$ws.Range("I2").Value = 0.5849329800180821
$ws.Range("J2").Value = 0.584932980018082
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1591403333333333
$ws.Range("N2").Value = 0.477421
$ws.Range("O2").Value = 0.01260326986877891
$ws.Range("P2").Value = 0.01260326986877891
$ws.Range("Q2").Value = 53.89673308946344
$ws.Range("R2").Value = 485.070597805171
$ws.Range("S2").Value = 0.00737206820231695
$ws.Range("T2").Value = 0.00737206820231695

# Row 3
$ws.Range("G3").Value = 338.6742503333333
$ws.Range("H3").Value = 1016.022751
$ws.Range("I3").Value = 0.5849329800180821
$ws.Range("J3").Value = 0.584932980018082
$ws.Range("O3").Value = 0.3005157372251983
$ws.Range("P3").Value = 0.3005157372251983
$ws.Range("Q3").Value = 1285.128117309695
$ws.Range("R3").Value = 11566.15305578725
$ws.Range("S3").Value = 0.1757815657174661
$ws.Range("T3").Value = 0.1757815657174661

# Row 4
$ws.Range("G4").Value = 338.6742503333333
$ws.Range("H4").Value = 1016.022751
$ws.Range("I4").Value = 0.5849329800180821
$ws.Range("J4").Value = 0.584932980018082
$ws.Range("M4").Value = 8.673183333333334
$ws.Range("N4").Value = 26.01955
$ws.Range("O4").Value = 0.6868809929060228
$ws.Range("P4").Value = 0.6868809929060229
$ws.Range("Q4").Value = 2937.383863420228
$ws.Range("R4").Value = 26436.45477078205
$ws.Range("S4").Value = 0.401779346098299
$ws.Range("T4").Value = 0.401779346098299

# Row 5
$ws.Range("I5").Value = 0.279688040971731
$ws.Range("J5").Value = 0.2796880409717309
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1591403333333333
$ws.Range("N5").Value = 0.477421
$ws.Range("O5").Value = 0.01260326986877891
$ws.Range("P5").Value = 0.01260326986877891
$ws.Range("Q5").Value = 25.77093822287522
$ws.Range("R5").Value = 231.938444005877
$ws.Range("S5").Value = 0.003524983859436818
$ws.Range("T5").Value = 0.003524983859436818

# Row 6
$ws.Range("I6").Value = 0.279688040971731
$ws.Range("J6").Value = 0.2796880409717309
$ws.Range("O6").Value = 0.3005157372251983
$ws.Range("P6").Value = 0.3005157372251983
$ws.Range("Q6").Value = 614.4891428705664
$ws.Range("R6").Value = 5530.402285835098
$ws.Range("S6").Value = 0.0840506578256912
$ws.Range("T6").Value = 0.0840506578256912

# Row 7
$ws.Range("I7").Value = 0.279688040971731
$ws.Range("J7").Value = 0.2796880409717309
$ws.Range("M7").Value = 8.673183333333334
$ws.Range("N7").Value = 26.01955
$ws.Range("O7").Value = 0.6868809929060228
$ws.Range("P7").Value = 0.6868809929060229
$ws.Range("Q7").Value = 1404.521827982039
$ws.Range("R7").Value = 12640.69645183835
$ws.Range("S7").Value = 0.192112399286603
$ws.Range("T7").Value = 0.192112399286603

# Row 8
$ws.Range("G8").Value = 77.79536166666666
$ws.Range("H8").Value = 233.386085
$ws.Range("I8").Value = 0.1343623634996766
$ws.Range("J8").Value = 0.1343623634996766
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1591403333333333
$ws.Range("N8").Value = 0.477421
$ws.Range("O8").Value = 0.01260326986877891
$ws.Range("P8").Value = 0.01260326986877891
$ws.Range("Q8").Value = 12.38037978742056
$ws.Range("R8").Value = 111.423418086785
$ws.Range("S8").Value = 0.001693405127393393
$ws.Range("T8").Value = 0.001693405127393393

# Row 9
$ws.Range("G9").Value = 77.79536166666666
$ws.Range("H9").Value = 233.386085
$ws.Range("I9").Value = 0.1343623634996766
$ws.Range("J9").Value = 0.1343623634996766
$ws.Range("O9").Value = 0.3005157372251983
$ws.Range("P9").Value = 0.3005157372251983
$ws.Range("Q9").Value = 295.2010865181211
$ws.Range("R9").Value = 2656.80977866309
$ws.Range("S9").Value = 0.04037800472242539
$ws.Range("T9").Value = 0.04037800472242539

# Row 10
$ws.Range("G10").Value = 77.79536166666666
$ws.Range("H10").Value = 233.386085
$ws.Range("I10").Value = 0.1343623634996766
$ws.Range("J10").Value = 0.1343623634996766
$ws.Range("M10").Value = 8.673183333333334
$ws.Range("N10").Value = 26.01955
$ws.Range("O10").Value = 0.6868809929060228
$ws.Range("P10").Value = 0.6868809929060229
$ws.Range("Q10").Value = 674.7334342179722
$ws.Range("R10").Value = 6072.600907961751
$ws.Range("S10").Value = 0.09229095364985783
$ws.Range("T10").Value = 0.09229095364985783

# Row 11
$ws.Range("G11").Value = 0.5886170000000001
$ws.Range("H11").Value = 1.765851
$ws.Range("I11").Value = 0.001016615510510267
$ws.Range("J11").Value = 0.001016615510510266
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1591403333333333
$ws.Range("N11").Value = 0.477421
$ws.Range("O11").Value = 0.01260326986877891
$ws.Range("P11").Value = 0.01260326986877891
$ws.Range("Q11").Value = 0.09367270558566668
$ws.Range("R11").Value = 0.843054350271
$ws.Range("S11").Value = [double]"1.281267963174733E-05"
$ws.Range("T11").Value = [double]"1.281267963174733E-05"

# Row 12
$ws.Range("G12").Value = 0.5886170000000001
$ws.Range("H12").Value = 1.765851
$ws.Range("I12").Value = 0.001016615510510267
$ws.Range("J12").Value = 0.001016615510510266
$ws.Range("O12").Value = 0.3005157372251983
$ws.Range("P12").Value = 0.3005157372251983
$ws.Range("Q12").Value = 2.233557042739334
$ws.Range("R12").Value = 20.102013384654
$ws.Range("S12").Value = 0.0003055089596155641
$ws.Range("T12").Value = 0.0003055089596155641

# Row 13
$ws.Range("G13").Value = 0.5886170000000001
$ws.Range("H13").Value = 1.765851
$ws.Range("I13").Value = 0.001016615510510267
$ws.Range("J13").Value = 0.001016615510510266
$ws.Range("M13").Value = 8.673183333333334
$ws.Range("N13").Value = 26.01955
$ws.Range("O13").Value = 0.6868809929060228
$ws.Range("P13").Value = 0.6868809929060229
$ws.Range("Q13").Value = 5.105183154116667
$ws.Range("R13").Value = 45.94664838705
$ws.Range("S13").Value = 0.0006982938712629552
$ws.Range("T13").Value = 0.0006982938712629551
